$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: updated publish date
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# A bare "true"/"false" typed into .Value is auto-coerced to a Boolean by
# this engine (matching Excel's literal-parsing rules), so force it to be
# stored as literal text via the classic apostrophe "text prefix" trick...
$c = $ws.Range("B17")
$c.Value = "'true"
# ...then strip the resulting quote-prefix formatting (which Excel applies
# automatically for apostrophe-forced text) by re-pasting the original
# cell's format from its still-untouched neighbour, so the style index is
# left exactly as it was.
$ws.Range("B16").Copy() | Out-Null
$c.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
